$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (Subjects / Tentative Exam Date shift right to H / I)
$ws.Columns.Item(7).Insert()

# New column G header/value: "Class" / "Std X" -- copy header style from neighboring header cell first
$ws.Range("H1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Class"
$ws.Range("G3").Value = "Std X"

# New column J header/value: "Reference" / "Rahul Dutta"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Reference"
$ws.Range("J3").Value = "Rahul Dutta"

# Move the legacy/threaded comment on the Subjects column from G3 to H3
$ws.Range("G3").Comment.Delete()
$ws.Range("H3").AddCommentThreaded("1) English`n2) French`n3) Physics`n4) Chemistry`n5) Mathematics`n6) Biology`n7) Computer`n")

# Column widths for the two new columns (best effort match of final layout)
$ws.Columns.Item(7).ColumnWidth = 4.25
$ws.Columns.Item(10).ColumnWidth = 9.92

# Match the saved selection state
$ws.Range("J4").Select()
